# "Comiteo template carga masiva"
# Bulk-upload template tweaks: clarify the two date-header columns with an
# explicit format hint, make room for a trailing marker column, and resize
# the existing columns so the wider headers/content fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clarify the expected date format directly in the column headers.
$ws.Range("D1").Value = "FECHA INICIO - dd/mm/yyyy hh:mm"
$ws.Range("E1").Value = "FECHA FIN -  dd/mm/yyyy hh:mm"

# Resize the data columns to fit the new header text / sample content.
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 34.5
$ws.Columns.Item(5).ColumnWidth = 34

# Touch a new trailing column (G) so it becomes part of the used range,
# extending the sheet's dimension/selection out to G1.
$ws.Range("G1").Font.Bold = $true

$ws.Range("G1").Select() | Out-Null
